$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("biochar_land")

# 1. Header label in G21 now reflects that the unit cost is per unit of INPUT
#    (non-energy costs are per unit input): "Unit cost  ($1975/kg)" -> "Unit cost  ($1975/kg input)"
$ws.Range("G21").Value = "Unit cost  ($1975/kg input)"

# 2. Copy the number formats from row 22 (Beef) down to the new row 27 (manure
#    total) so it matches the existing per-animal rows' styles.
$ws.Range("B22:G22").Copy()
$ws.Range("B27:G27").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# 3. Label A27 "manure" (was a blank placeholder row before). Stash/restore its
#    existing style first, since writing .Value resets the cell's quote-prefix
#    (text) formatting that the original blank cell carried.
$ws.Range("A27").Copy($ws.Range("Z1"))
$ws.Range("A27").Value = "manure"
$ws.Range("Z1").Copy()
$ws.Range("A27").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("Z1").Clear()

# 4. Populate the new "manure" total row (B27:G27), mirroring rows 22-26.
#    Unlike the per-animal rows, F27 (production estimate) isn't scaled by a
#    per-animal bio-oil yield factor -- it's the aggregate across all manure.
$ws.Range("B27").Formula = "=(A`$8+A`$9)*A`$10"
$ws.Range("C27").Formula = "=0.09*B27"
$ws.Range("D27").Formula = "=-PV(A`$11,A`$5,C27)+B27"
$ws.Range("E27").Formula = "=-PMT(A`$11,A`$5,D27)"
$ws.Range("F27").Formula = "=A`$3*A`$4*A`$6*A`$7*A`$18"
$ws.Range("F27").Style = "Normal"
$ws.Range("G27").Formula = "=E27/F27"

# 5. Leave the selection where the author last clicked when saving.
$ws.Range("G23").Select()
